$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: push the existing "sum" row (7) and trailing blank row (8) down to rows 9/10,
#        carrying their formatting with them. Row 8 originally only has a formatted H8 cell, so
#        only move that one cell (otherwise Copy materialises empty <c> tags for the whole row). ---
$ws.Range("A7:I7").Copy($ws.Range("A9:I9"))
$ws.Range("H8").Copy($ws.Range("H10"))
$ws.Range("H8").ClearContents()

# --- 2. Seed rows 7 and 8 with the same formatting as row 6 (a normal data row), then clear
#        their content so we can fill in the new parts cleanly. ---
$ws.Range("A6:I6").Copy($ws.Range("A7:I7"))
$ws.Range("A6:I6").Copy($ws.Range("A8:I8"))
$ws.Range("A7:I8").ClearContents()

# --- 3. New part: SparkFun Carrier Board (row 7) ---
$ws.Range("B7").Value = "SparkFun Carrier Board"
$ws.Range("G7").Value = "https://www.digikey.ch/de/products/detail/sparkfun-electronics/16885/13282886?s=N4IgTCBcDaICIFEBqBaAjANgBxYKwgF0BfIA"
$ws.Range("E7").Value = "1568-16885-ND"
$ws.Range("F7").ClearFormats()
$ws.Range("F7").Value = "DEV-16885"
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "SparkFun"
$ws.Range("D7").Value = "Digikey"
$ws.Range("H7").Value = 20.5

# --- 4. New part: SparkFun ESP32 Processor (row 8) ---
$ws.Range("B8").Value = "SparkFun ESP32 Processor"
$ws.Range("F8").ClearFormats()
$ws.Range("F8").Value = "WRL-16781"
$ws.Range("G8").Value = "https://www.digikey.ch/de/products/detail/sparkfun-electronics/16781/13282892?s=N4IgTCBcDaIOoCUAyBaAjANgOwA40gF0BfIA"
$ws.Range("E8").Value = "1568-16781-ND"
$ws.Range("A8").Value = 1
$ws.Range("C8").Value = "SparkFun"
$ws.Range("D8").Value = "Digikey"
$ws.Range("H8").Value = 20.6

# --- 5. Begründung text is filled in last for both new rows. ---
$ws.Range("I7").Value = "Carrier Board um den ESP32 Prozessor zu halten"
$ws.Range("I8").Value = "ESP32 Prozessor um alles zu rechnen"

# --- 6. Fix up the (now-moved) totals row 9: recompute the sum over the extended range. ---
$ws.Range("A9:I9").ClearContents()
$ws.Range("H9").Formula = "=SUM(H2:H8)"

# --- 7. Row 10 is the trailing blank row (only H10 carries formatting) - already moved/cleared. ---

# --- 8. Restore the selection like the saved file shows. ---
$ws.Range("I8").Select()
